$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (column D stores
# thousands-dot-formatted price strings, not real numbers) -- set the
# format per-cell since multi-area Range() selections only apply to the
# first area in this COM layer.
$textCells = @("D5", "D6", "D8", "D11", "D15", "D16", "D18", "D19", "D25", "D27", "D32", "D37", "D39", "D40", "D41", "D47", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.936.60'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '1.646.37'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '213.53'
$ws.Range('D6').Value = '0.527'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '23.43'
$ws.Range('E8').Value = '  +2.60%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '0.0871'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').Value = '1.881.05'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').Value = '1.648.89'
$ws.Range('E13').Value = '  +1.86%  '
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').Value = '65.65'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '27.940.55'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '232.59'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '7.66'
$ws.Range('E19').Value = '  +1.88%  '
$ws.Range('D20').Value = '0.0₃0722'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  +4.61%  '
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('E24').Value = '  +4.18%  '
$ws.Range('D25').Value = '152.47'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = '15.74'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').Value = '3.36'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').Value = '1.449.48'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').Value = '0.888'
$ws.Range('E37').Value = '  +2.73%  '
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').Value = '0.562'
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').Value = '0.920'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').Value = '69.44'
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E46').Value = '  +5.31%  '
$ws.Range('D47').Value = '5.35'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = '1.788.92'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').Value = '89.04'
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('E51').Value = '  +0.70%  '
